# Improve app to work with Tables UI tests
# - add three new survey fields (db, services, all_in_one) right before
#   the existing "time" (Test step / Time of step) rows
# - update the saved window position / active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert 3 fresh rows above the current row 10 ("time" field); Excel's
# Insert() shifts the existing rows (and their formatting) down and the
# new rows inherit the formatting of the row above (row 9), matching
# s="4" (col A) / s="5" (cols B-E) used by the rest of the "text" fields.
$ws.Rows("10:12").Insert()

# New field: db / "Database Used"
$ws.Range("A10").Value = "text"
$ws.Range("C10").Value = "db"
$ws.Range("D10").Value = "Database Used"

# New field: services / "Services Used"
$ws.Range("A11").Value = "text"
$ws.Range("C11").Value = "services"
$ws.Range("D11").Value = "Services Used"

# New field: all_in_one / "All In One APK Used"
$ws.Range("A12").Value = "text"
$ws.Range("C12").Value = "all_in_one"
$ws.Range("D12").Value = "All In One APK Used"

# Match the saved selection/active cell in the authored workbook.
$ws.Range("E12").Select()

# Best-effort: restore the saved window position recorded in the
# workbook view (cosmetic - where the Excel window sat on screen).
$win = $excel.ActiveWindow
$win.Left = 25280
$win.Top = 5960
